$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.618.34"
$ws.Range("E2").Value = "  +3.92%  "

$ws.Range("D3").Value = "1.798.68"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5306"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.94%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3768"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07529"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.57"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.122"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.21"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.58%  "

$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.207"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.478"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.87%  "

$ws.Range("D16").Value = "1.796.35"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001071"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06464"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.924"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("D23").Value = "28.640.04"
$ws.Range("E23").Value = "  +3.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.18"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.02"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.55"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.378"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.45%  "

$ws.Range("D29").Value = "2.001.67"
$ws.Range("E29").Value = "  +0.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.61"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.127"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1025"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.724"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.669"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2309"
$ws.Range("D35").ClearFormats()

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06539"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.80%  "

$ws.Range("E37").Value = "  +2.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.830"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.069"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.50"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6319"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.203"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.99%  "

$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.393"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.50"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5932"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.670"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.93"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.47%  "

$ws.Range("E49").Value = "  +3.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.170"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06931"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.67%  "
